$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Macbook retina Python 3.4 64bity" (numpy) benchmark block ---

# Section header, row 11 (mirrors A7's section-header pattern)
$ws.Range("A11").Value = "Macbook retina Python 3.4 64bity"

# Raw timings, row 12 (mirrors row 8's data pattern)
$ws.Range("F12").Value = 7.32
$ws.Range("G12").Value = 4.01

# Relative-speed formulas, row 13 (mirrors row 9's shared-formula pattern)
$ws.Range("F13:G13").Formula = "=F12/`$F`$3"

# Copy the percent number format from the existing F9:G9 formula cells
# onto the new F13:G13 cells without disturbing their formulas/values.
$ws.Range("F9:G9").Copy()
$ws.Range("F13:G13").PasteSpecial(-4122)

# Move the active selection to G22, matching the final saved cursor position
$ws.Range("G22").Select() | Out-Null
